$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date for rows 2-29 from 45559 to 45560
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 3).Value = 45560
}

# Row 29 loses its explicit custom row height (reverts to default/auto height)
$ws.Rows.Item(29).AutoFit()

# Delete row 30 entirely (last data row removed)
$ws.Rows.Item(30).Delete()
